# regen sval data to filter save games
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New values for columns B, C, D, E, G across rows 2-14.
# (Column A = date labels, F = win flag; both unchanged.)
$data = @{
    2  = @{ B = 1.505614041169197;  C = 1.65323645889881;   D = 0.1529057820181812;  E = 0.4998867070740569; G = 3.811642989160245 }
    3  = @{ B = 3.182878228561681;  C = 1.65323645889881;   D = 0.7127328510149897;  E = 0.4998867070740569; G = 6.048734245549538 }
    4  = @{ B = 3.182878228561681;  C = 1.65323645889881;   D = 0.7127328510149897;  E = 0.4998867070740569; G = 6.048734245549538 }
    5  = @{ B = 3.182878228561681;  C = 1.65323645889881;   D = 0.7127328510149897;  E = 0.4998867070740569; G = 6.048734245549538 }
    6  = @{ B = 0.7287194209349384; C = 0.05231270169004087; D = 0.1529057820181812; E = 0.4998867070740569; G = 1.433824611717217 }
    7  = @{ B = 3.182878228561681;  C = 1.65323645889881;   D = 16.98373111632243;  E = 0.4998867070740569; G = 22.31973251085698 }
    8  = @{ B = 3.182878228561681;  C = 1.65323645889881;   D = 3.082599426703578;  E = 0.4998867070740569; G = 8.418600821238126 }
    9  = @{ B = 3.182878228561681;  C = 1.65323645889881;   D = 0.7127328510149897;  E = 0.4998867070740569; G = 6.048734245549538 }
    10 = @{ B = 0.7287194209349384; C = 1.65323645889881;   D = 0.7127328510149897;  E = 0.4998867070740569; G = 3.594575437922795 }
    11 = @{ B = 0.06328177979961902; C = 0.05231270169004087; D = 0.1529057820181812; E = 0.4998867070740569; G = 0.768386970581898 }
    12 = @{ B = 0.1554434735375247; C = 0.3375848360084654;  D = 0.1529057820181812; E = 0.4998867070740569; G = 1.145820798638228 }
    13 = @{ B = 1.505614041169197;  C = 1.65323645889881;   D = 3.082599426703578;  E = 0.4998867070740569; G = 6.741336633845642 }
    14 = @{ B = 0.7287194209349384; C = 0.3375848360084654;  D = 0.1529057820181812; E = 0.4998867070740569; G = 1.719096746035642 }
}

foreach ($row in $data.Keys) {
    $vals = $data[$row]
    $ws.Range("B$row").Value = $vals.B
    $ws.Range("C$row").Value = $vals.C
    $ws.Range("D$row").Value = $vals.D
    $ws.Range("E$row").Value = $vals.E
    $ws.Range("G$row").Value = $vals.G
}
